# Rename "Sheet1" to "Data"
$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "Data"

# On the "Legend" sheet, insert a new header row above the existing
# ID/label pairs and turn the range into an Excel Table ("Table1") with
# generic "Column1"/"Column2" headers.
$wsLegend = $wb.Worksheets.Item("Legend")
$wsLegend.Rows.Item(1).Insert()
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

$tableRange = $wsLegend.Range("A1:B7")
$lo = $wsLegend.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# Select the whole table on the Legend sheet (matches the saved selection
# state), then restore the originally-active "Data" sheet/tab.
$wsLegend.Activate()
$tableRange.Select() | Out-Null
$wsData.Activate()
